$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Replace the "Participants" tab's Neo4j (dbExcel) query with the updated
# Cypher statement. The old query text is removed from the shared-string
# table automatically (it becomes unreferenced) and the new one is appended,
# which is exactly what shifts the other rows' StatQuery/FilesQuery shared
# string indices in the underlying XML - no other cells need to be touched.
$newParticipantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['TXT']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$ws.Range("B2").Value = $newParticipantsQuery

# Move the active selection from C5 to B5, scrolled so row 4 is at the top
# of the view (matches the updated sheetView/selection in the workbook).
$win = $excel.ActiveWindow
[void]$ws.Range("B5").Select()
$win.ScrollRow = 4
$win.ScrollColumn = 1
